$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the current row 173, which
# pushes the existing rows 173-179 down to 174-180 (dimension grows to
# A1:R180). The newly inserted row carries the latest observation.
$ws.Rows("173:173").Insert()

$ws.Cells.Item(173, 1).Value = 9
$ws.Cells.Item(173, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(173, 3).Value = "Metropolitana"
$ws.Cells.Item(173, 4).Value = 44509
$ws.Cells.Item(173, 5).Value = 13
$ws.Cells.Item(173, 6).Value = 300000001
$ws.Cells.Item(173, 7).Value = "Rabanito"
$ws.Cells.Item(173, 8).Value = "Sin especificar"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 8800
$ws.Cells.Item(173, 11).Value = 2500
$ws.Cells.Item(173, 12).Value = 3000
$ws.Cells.Item(173, 13).Value = 2750
$ws.Cells.Item(173, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(173, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(173, 16).Value = 28
$ws.Cells.Item(173, 17).Value = 100
$ws.Cells.Item(173, 18).Value = "Hortaliza"
